# Update ESD excel file
# - Remove the two "Number of fetuses" rows (previously rows 33 and 34,
#   concept ids 40758410 and 3002549) which were the last rows of the
#   "observationConceptIds" category.
# - Re-categorize the remaining three "observationConceptIds" rows
#   (concept ids 3011536, 3026070, 3024261) to "estDeliveryConceptIds".
# - Table1 / sheet dimension / selection shrink accordingly (A1:D34 -> A1:D32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two trailing rows for "Number of fetuses" / "Number of fetuses by US".
$ws.Rows("33:34").Delete()

# Re-point the remaining "observationConceptIds" rows to "estDeliveryConceptIds".
$ws.Range("D30:D32").Value2 = "estDeliveryConceptIds"

# Move the active selection to D32, matching the final saved selection.
$ws.Range("D32").Select()
